$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.921.26"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "1.810.50"
$ws.Range("E3").Value = "  +1.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("E7").Value = "  -2.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3909"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09945"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +26.92%  "

$ws.Range("E10").Value = "  +0.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.396"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.92%  "

$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.000"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "

$ws.Range("D15").Value = "1.808.67"
$ws.Range("E15").Value = "  +1.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.250"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001146"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06645"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.920"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").Value = "27.986.80"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.262"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.51%  "

$ws.Range("D27").Value = "2.020.35"
$ws.Range("E27").Value = "  +1.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.392"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1060"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.027"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.554"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.591"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06702"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.90%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02323"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.93%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.865"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.06%  "

$ws.Range("E38").Value = "  +0.73%  "

$ws.Range("E39").Value = "  -1.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6169"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.168"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5882"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.682"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.278"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.927"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.175"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06767"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
